$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 977.1053000000001
$ws.Range("I2").Value = 718.61536
$ws.Range("J2").Value = 1537.1666
$ws.Range("K2").Value = 718.61536
$ws.Range("L2").Value = 1537.1666
$ws.Range("M2").Value = -605.61536
$ws.Range("N2").Value = -1763.1666
$ws.Range("H76").Value = 3316.9211
$ws.Range("I76").Value = 2455.7407
$ws.Range("J76").Value = 5430.727
$ws.Range("K76").Value = 2455.7407
$ws.Range("L76").Value = 5430.727
$ws.Range("M76").Value = -2140.7407
$ws.Range("N76").Value = -6060.727
$ws.Range("H79").Value = 3316.9211
$ws.Range("I79").Value = 2455.7407
$ws.Range("J79").Value = 5430.727
$ws.Range("K79").Value = 2455.7407
$ws.Range("L79").Value = 5430.727
$ws.Range("M79").Value = -1363.7407
$ws.Range("N79").Value = -7614.727
$ws.Range("H98").Value = 3639.7058
$ws.Range("I98").Value = 2304.6875
$ws.Range("J98").Value = 25000
$ws.Range("K98").Value = 2304.6875
$ws.Range("L98").Value = 25000
$ws.Range("M98").Value = -806.6875
$ws.Range("N98").Value = -27996
$ws.Range("H122").Value = 3639.7058
$ws.Range("I122").Value = 2304.6875
$ws.Range("J122").Value = 25000
$ws.Range("K122").Value = 6914.0625
$ws.Range("L122").Value = 75000
$ws.Range("M122").Value = -4464.0625
$ws.Range("N122").Value = -79900
$ws.Range("H132").Value = 22399.941
$ws.Range("I132").Value = 22399.941
$ws.Range("K132").Value = 67199.823
$ws.Range("M132").Value = -64669.823
$ws.Range("H135").Value = 2774.476
$ws.Range("I135").Value = 2324.2666
$ws.Range("K135").Value = 20918.3994
$ws.Range("M135").Value = -18383.3994
$ws.Range("H138").Value = 58633.277
$ws.Range("I138").Value = 3213.3333
$ws.Range("J138").Value = 86343.25
$ws.Range("K138").Value = 9639.999899999999
$ws.Range("L138").Value = 259029.75
$ws.Range("M138").Value = -4499.999899999999
$ws.Range("N138").Value = -269309.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23689.744
$ws.Range("I32").Value = 25733.162
$ws.Range("K32").Value = 25733.162
$ws.Range("M32").Value = -25446.162
$ws.Range("H45").Value = 3309.8823
$ws.Range("I45").Value = 1671.125
$ws.Range("K45").Value = 1671.125
$ws.Range("M45").Value = -1294.125
$ws.Range("H61").Value = 7221.647
$ws.Range("I61").Value = 996.7143
$ws.Range("K61").Value = 996.7143
$ws.Range("M61").Value = -784.7143
$ws.Range("H63").Value = 2959.8667
$ws.Range("I63").Value = 2812.8333
$ws.Range("K63").Value = 2812.8333
$ws.Range("M63").Value = -2126.8333
$ws.Range("H66").Value = 2959.8667
$ws.Range("I66").Value = 2812.8333
$ws.Range("K66").Value = 14064.1665
$ws.Range("M66").Value = -10632.1665
$ws.Range("H74").Value = 360091.25
$ws.Range("I74").Value = 750784.6
$ws.Range("K74").Value = 750784.6
$ws.Range("M74").Value = -749910.6
$ws.Range("H77").Value = 360091.25
$ws.Range("I77").Value = 750784.6
$ws.Range("K77").Value = 3753923
$ws.Range("M77").Value = -3749555
$ws.Range("H132").Value = 1572.3438
$ws.Range("I132").Value = 1111.52
$ws.Range("J132").Value = 3218.1428
$ws.Range("K132").Value = 3334.56
$ws.Range("L132").Value = 9654.428400000001
$ws.Range("M132").Value = -804.5599999999999
$ws.Range("N132").Value = -14714.4284
$ws.Range("H136").Value = 7221.647
$ws.Range("I136").Value = 996.7143
$ws.Range("K136").Value = 2990.1429
$ws.Range("M136").Value = -440.1428999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 627.3333
$ws.Range("I22").Value = 404.83334
$ws.Range("J22").Value = 849.8333
$ws.Range("K22").Value = 404.83334
$ws.Range("L22").Value = 849.8333
$ws.Range("M22").Value = -54.83334000000002
$ws.Range("N22").Value = -1549.8333
$ws.Range("H31").Value = 14289467
$ws.Range("I31").Value = 33334134
$ws.Range("J31").Value = 5967.5
$ws.Range("K31").Value = 33334134
$ws.Range("L31").Value = 5967.5
$ws.Range("M31").Value = -33333839
$ws.Range("N31").Value = -6557.5
$ws.Range("H34").Value = 14289467
$ws.Range("I34").Value = 33334134
$ws.Range("J34").Value = 5967.5
$ws.Range("K34").Value = 33334134
$ws.Range("L34").Value = 5967.5
$ws.Range("M34").Value = -33333932
$ws.Range("N34").Value = -6371.5
$ws.Range("H132").Value = 113845.555
$ws.Range("I132").Value = 201172
$ws.Range("K132").Value = 603516
$ws.Range("M132").Value = -600986
$ws.Range("H134").Value = 3001.1304
$ws.Range("I134").Value = 2823.0527
$ws.Range("K134").Value = 8469.158100000001
$ws.Range("M134").Value = -5934.158100000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 44450.75
$ws.Range("J37").Value = 44450.75
$ws.Range("L37").Value = 133352.25
$ws.Range("N37").Value = -133576.25
$ws.Range("H107").Value = 2645.9333
$ws.Range("I107").Value = 6258.2
$ws.Range("J107").Value = 839.8
$ws.Range("K107").Value = 18774.6
$ws.Range("L107").Value = 2519.4
$ws.Range("M107").Value = -16854.6
$ws.Range("N107").Value = -6359.4

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 31250
$ws.Range("J74").Value = 31250
$ws.Range("L74").Value = 31250
$ws.Range("N74").Value = -33122
$ws.Range("H77").Value = 31250
$ws.Range("J77").Value = 31250
$ws.Range("L77").Value = 93750
$ws.Range("N77").Value = -103110
$ws.Range("H80").Value = 9750.237999999999
$ws.Range("I80").Value = 5609.7856
$ws.Range("J80").Value = 18031.143
$ws.Range("K80").Value = 5609.7856
$ws.Range("L80").Value = 18031.143
$ws.Range("M80").Value = -4611.7856
$ws.Range("N80").Value = -20027.143
$ws.Range("H83").Value = 9750.237999999999
$ws.Range("I83").Value = 5609.7856
$ws.Range("J83").Value = 18031.143
$ws.Range("K83").Value = 28048.928
$ws.Range("L83").Value = 90155.715
$ws.Range("M83").Value = -23056.928
$ws.Range("N83").Value = -100139.715
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H95").Value = 69661.336
$ws.Range("J95").Value = 69661.336
$ws.Range("L95").Value = 69661.336
$ws.Range("N95").Value = -75153.336
$ws.Range("H126").Value = 3166.6667
$ws.Range("I126").Value = 1750
$ws.Range("K126").Value = 5250
$ws.Range("M126").Value = -2780
$ws.Range("H132").Value = 2111.8538
$ws.Range("I132").Value = 2014.6177
$ws.Range("J132").Value = 2584.1428
$ws.Range("K132").Value = 6043.8531
$ws.Range("L132").Value = 7752.428400000001
$ws.Range("M132").Value = -3513.8531
$ws.Range("N132").Value = -12812.4284

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1640.2222
$ws.Range("I55").Value = 396
$ws.Range("K55").Value = 396
$ws.Range("M55").Value = -223
$ws.Range("H122").Value = 2287.125
$ws.Range("I122").Value = 2372.9333
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 7118.7999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -4668.7999
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 2949.75
$ws.Range("I132").Value = 2718
$ws.Range("J132").Value = 3830.4
$ws.Range("K132").Value = 8154
$ws.Range("L132").Value = 11491.2
$ws.Range("M132").Value = -5624
$ws.Range("N132").Value = -16551.2

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 60448.56
$ws.Range("I122").Value = 71163.48
$ws.Range("K122").Value = 213490.44
$ws.Range("M122").Value = -211040.44
$ws.Range("H132").Value = 26037.535
$ws.Range("I132").Value = 27617.424
$ws.Range("J132").Value = 5499
$ws.Range("K132").Value = 82852.272
$ws.Range("L132").Value = 16497
$ws.Range("M132").Value = -80322.272
$ws.Range("N132").Value = -21557
$ws.Range("H136").Value = 16389.129
$ws.Range("I136").Value = 21395.785
$ws.Range("K136").Value = 64187.355
$ws.Range("M136").Value = -61637.355
